$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = "NOTE: populate JP1-4 depending on whether you want the inputs to be AC or DC coupled. For DC coupling, wire links or zero"
$ws.Range("E16").Value = "ohm links will work. For AC coupling, use a 100nF ceramic capacitor – if you plan to exploit mixer feedback then AC coupling"
$ws.Range("E17").Value = "is recommended in order to prevent the op-amps latching up."

$ws.Range("A1").Select()
$ws.Range("E20").Select()
